# Apply updated crypto price/volume figures to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Price" column (D): values look numeric (contain dots as thousand
# separators), but must stay as literal text, exactly as authored.
# Temporarily force a Text number format so Excel does not coerce the
# string into a float (which would silently drop formatting such as
# "212.89" groups or trailing zeros like "0.510"), then restore the
# original (default/"Normal") style so no stray formatting is left behind.
$priceUpdates = @{
    "D2" = "26.261.50"
    "D3" = "1.593.15"
    "D5" = "212.89"
    "D11" = "0.0851"
    "D12" = "1.816.30"
    "D13" = "1.581.88"
    "D15" = "0.510"
    "D16" = "63.90"
    "D17" = "26.257.99"
    "D20" = "7.35"
    "D25" = "144.68"
    "D27" = "6.97"
    "D29" = "15.13"
    "D33" = "1.414.79"
    "D34" = "2.97"
    "D40" = "5.82"
    "D42" = "0.970"
    "D45" = "1.729.05"
    "D46" = "60.98"
    "D47" = "86.82"
    "D49" = "0.0502"
    "D50" = "0.0959"
    "D51" = "0.999"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceUpdates[$addr]
    $ws.Range($addr).Style = "Normal"
}

# --- "Volume(1h)" column (E): percentage text like "  +0.04%  " is never
# numeric-looking to Excel, so a plain Value assignment is sufficient.
$volumeUpdates = @{
    "E3" = "  +0.04%  "
    "E4" = "  +0.01%  "
    "E5" = "  +0.59%  "
    "E6" = "  -0.78%  "
    "E7" = "  +0.03%  "
    "E8" = "  -0.57%  "
    "E9" = "  -0.59%  "
    "E10" = "  -2.45%  "
    "E11" = "  +0.50%  "
    "E12" = "  -0.01%  "
    "E13" = "  -0.77%  "
    "E14" = "  -1.64%  "
    "E15" = "  -2.76%  "
    "E16" = "  -1.32%  "
    "E19" = "  +1.04%  "
    "E20" = "  -2.26%  "
    "E21" = "  +0.01%  "
    "E22" = "  -0.29%  "
    "E23" = "  +0.16%  "
    "E24" = "  -2.59%  "
    "E25" = "  +0.50%  "
    "E26" = "  +0.01%  "
    "E27" = "  -1.69%  "
    "E28" = "  -0.99%  "
    "E29" = "  -0.77%  "
    "E30" = "  -2.69%  "
    "E31" = "  +0.31%  "
    "E32" = "  -0.69%  "
    "E33" = "  +5.77%  "
    "E35" = "  -0.54%  "
    "E36" = "  -1.18%  "
    "E37" = "  -3.14%  "
    "E38" = "  -1.06%  "
    "E39" = "  +0.39%  "
    "E41" = "  +0.05%  "
    "E42" = "  -5.61%  "
    "E43" = "  +0.35%  "
    "E44" = "  -0.32%  "
    "E45" = "  -0.04%  "
    "E46" = "  -1.59%  "
    "E47" = "  -1.51%  "
    "E48" = "  -0.38%  "
    "E49" = "  -0.61%  "
    "E50" = "  -2.81%  "
    "E51" = "  +0.01%  "
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
